$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the "Running?" flag (Yes) on case6_1 (row13) and case6_3 (row15) ---
$ws.Range("G13").Value = ""
$ws.Range("G15").Value = ""

# --- Remove the old row 63 cleanly so it can be rebuilt as row 64 with new data ---
$ws.Rows.Item(63).Delete() | Out-Null

# --- Row 64: case6_11, now continuing with subcase 11, updated remark/result text ---
$ws.Range("B64").Value = "my laptop"
$ws.Range("C64").Value = "case6_11"
$ws.Range("D64").Value = 6
$ws.Range("E64").Value = 11
$ws.Range("G64").Font.Color = 255
$ws.Range("G64").HorizontalAlignment = -4108
$ws.Range("H64").HorizontalAlignment = -4108
$ws.Range("H64").Value = "'""relu"""
$ws.Range("I64").Value = "CNN: added one line of dense at the final"
$ws.Range("J64").WrapText = $true
$ws.Range("J64").Font.Color = 255
$ws.Range("J64").Value = "good results, saved"

# --- Row 65: case6_12 ---
$ws.Range("B65").Value = 811228
$ws.Range("C65").Value = "case6_12"
$ws.Range("D65").Value = 6
$ws.Range("E65").Value = 12
$ws.Range("G65").Font.Color = 255
$ws.Range("G65").HorizontalAlignment = -4108
$ws.Range("G65").Value = "Yes"
$ws.Range("H65").HorizontalAlignment = -4108
$ws.Range("H65").Value = "'""relu"""
$ws.Range("I65").Value = "CNN: added one line of dense at the final"

# --- Row 66: case6_13 ---
$ws.Range("B66").Value = 811229
$ws.Range("C66").Value = "case6_13"
$ws.Range("D66").Value = 6
$ws.Range("E66").Value = 13
$ws.Range("G66").Font.Color = 255
$ws.Range("G66").HorizontalAlignment = -4108
$ws.Range("G66").Value = "Yes"
$ws.Range("H66").HorizontalAlignment = -4108
$ws.Range("H66").Value = "'""relu"""
$ws.Range("I66").Value = "CNN: added one line of dense at the final"

# --- Row 67: case6_14 ---
$ws.Range("B67").Value = 811233
$ws.Range("C67").Value = "case6_14"
$ws.Range("D67").Value = 6
$ws.Range("E67").Value = 14
$ws.Range("G67").Font.Color = 255
$ws.Range("G67").HorizontalAlignment = -4108
$ws.Range("G67").Value = "Yes"
$ws.Range("H67").HorizontalAlignment = -4108
$ws.Range("H67").Value = "'""relu"""
$ws.Range("I67").Value = "CNN: added one line of dense at the final"

# --- Row 68: case6_15 ---
$ws.Range("B68").Value = 811232
$ws.Range("C68").Value = "case6_15"
$ws.Range("D68").Value = 6
$ws.Range("E68").Value = 15
$ws.Range("G68").Font.Color = 255
$ws.Range("G68").HorizontalAlignment = -4108
$ws.Range("G68").Value = "Yes"
$ws.Range("H68").HorizontalAlignment = -4108
$ws.Range("H68").Value = "'""relu"""
$ws.Range("I68").Value = "CNN: added one line of dense at the final"

# --- Row 69: blank new row, just carries H style forward ---
$ws.Range("H69").Value = "'x"
$ws.Range("H69").HorizontalAlignment = -4108
$ws.Range("H69").Value = ""

# --- Restore the selected cell as in the final workbook ---
$ws.Range("G54").Select() | Out-Null
